$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data had several match rows shuffled out of chronological/
# listing order. This script fixes the ordering by rotating the match
# details (columns F:V -- everything except the Indice/pais/torneio/
# temporada/data_partida columns A:E, which stay put) among the affected
# rows, and appends one missing match row (113) that had been dropped.

# NOTE: this PowerShell host does not support named (-Param value) argument
# binding on user-defined functions, only positional, so Rotate-Rows takes
# its row list positionally.
function Rotate-Rows {
    param($Rows)
    # Rotate F:V content through the given rows: new(Rows[0]) = old(Rows[1]),
    # new(Rows[1]) = old(Rows[2]), ..., new(Rows[-1]) = old(Rows[0]).
    $snapshots = @()
    foreach ($r in $Rows) {
        $snapshots += ,($ws.Range("F${r}:V${r}").Value2)
    }
    $n = $Rows.Length
    for ($i = 0; $i -lt $n; $i++) {
        $targetRow = $Rows[$i]
        $sourceData = $snapshots[($i + 1) % $n]
        $ws.Range("F${targetRow}:V${targetRow}").Value = $sourceData
    }
}

# Rows 69 -> 70 -> 72 -> 69 (3-way rotation)
Rotate-Rows @(69, 70, 72)

# Rows 83 -> 84 -> 85 -> 83 (3-way rotation)
Rotate-Rows @(83, 84, 85)

# Rows 93 <-> 95 (swap)
Rotate-Rows @(93, 95)

# Rows 103 <-> 104 (swap)
Rotate-Rows @(103, 104)

# Rows 105 <-> 106 (swap)
Rotate-Rows @(105, 106)

# Append the missing match as new row 113, copying the formatting from the
# last existing row (112) so the bold/centered index column and the
# datetime-formatted date column keep their look.
$ws.Range("A112:V112").Copy()
$ws.Range("A113").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A113").Value = 112
$ws.Range("B113").Value = "poland"
$ws.Range("C113").Value = "iii-liga-group-iii"
$ws.Range("D113").Value = "2023-2024"
$ws.Range("E113").Value = 45235.5
$ws.Range("F113").Value = "Gornik Zabrze II"
$ws.Range("G113").Value = 1
$ws.Range("H113").Value = "Bielsko-Biala"
$ws.Range("I113").Value = 3
$ws.Range("J113").Value = 3.36
$ws.Range("K113").Value = "04/11/2023 00:14"
$ws.Range("L113").Value = 3.55
$ws.Range("M113").Value = "05/11/2023 09:51"
$ws.Range("N113").Value = 3.67
$ws.Range("O113").Value = "04/11/2023 00:14"
$ws.Range("P113").Value = 3.86
$ws.Range("Q113").Value = "05/11/2023 11:01"
$ws.Range("R113").Value = 1.72
$ws.Range("S113").Value = "04/11/2023 00:14"
$ws.Range("T113").Value = 1.76
$ws.Range("U113").Value = "05/11/2023 09:51"
$ws.Range("V113").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/gornik-zabrze-rekord-bielsko-biala/ruEq2wIF/"
